# Auto-generated Excel COM-interop script
# Applies numeric odds updates to Sheet1 as described in the commit diff
# Commit message: Atualizando o arquivo XLSX

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 1.18
$ws.Range("K2").Value = 4.5
$ws.Range("L2").Value = 1.83
$ws.Range("M2").Value = 1.83
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 1.88
$ws.Range("Q2").Value = 1.93

# Row 3
$ws.Range("AD3").Value = 451
$ws.Range("I3").Value = 3.75
$ws.Range("J3").Value = 1.08
$ws.Range("K3").Value = 8
$ws.Range("N3").Value = 2.3
$ws.Range("O3").Value = 1.6
$ws.Range("P3").Value = 1.5
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 1.75
$ws.Range("Z3").Value = 7.5

# Row 4
$ws.Range("AH4").Value = 29
$ws.Range("G4").Value = 2.7
$ws.Range("I4").Value = 2.7

# Row 6
$ws.Range("AA6").Value = 7
$ws.Range("AD6").Value = 201
$ws.Range("AE6").Value = 8
$ws.Range("AF6").Value = 8.5
$ws.Range("AJ6").Value = 23
$ws.Range("H6").Value = 3.75
$ws.Range("I6").Value = 1.7
$ws.Range("K6").Value = 12
$ws.Range("N6").Value = 1.75
$ws.Range("O6").Value = 2.05
$ws.Range("R6").Value = 1.75
$ws.Range("S6").Value = 2
$ws.Range("Y6").Value = 41
$ws.Range("Z6").Value = 12

# Row 7
$ws.Range("AB7").Value = 15
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 51
$ws.Range("AI7").Value = 34
$ws.Range("AJ7").Value = 34
$ws.Range("G7").Value = 1.7
$ws.Range("L7").Value = 1.2
$ws.Range("M7").Value = 4.33
$ws.Range("N7").Value = 1.67
$ws.Range("O7").Value = 2.15
$ws.Range("R7").Value = 1.67
$ws.Range("S7").Value = 2.1
$ws.Range("T7").Value = 9
$ws.Range("U7").Value = 9
$ws.Range("W7").Value = 13
$ws.Range("Z7").Value = 15

# Row 8
$ws.Range("AB8").Value = 21
$ws.Range("AC8").Value = 81
$ws.Range("AE8").Value = 9.5
$ws.Range("AF8").Value = 21
$ws.Range("AH8").Value = 51
$ws.Range("AJ8").Value = 51
$ws.Range("G8").Value = 1.91
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 4.75
$ws.Range("J8").Value = 1.1
$ws.Range("K8").Value = 7
$ws.Range("L8").Value = 1.44
$ws.Range("M8").Value = 2.63
$ws.Range("N8").Value = 2.5
$ws.Range("O8").Value = 1.5
$ws.Range("P8").Value = 1.57
$ws.Range("Q8").Value = 2.25
$ws.Range("R8").Value = 2.2
$ws.Range("S8").Value = 1.62
$ws.Range("T8").Value = 5.5
$ws.Range("U8").Value = 7.5
$ws.Range("V8").Value = 9.5
$ws.Range("W8").Value = 15
$ws.Range("X8").Value = 19
$ws.Range("Y8").Value = 41
$ws.Range("Z8").Value = 6

# Row 9
$ws.Range("AA9").Value = 6
$ws.Range("AJ9").Value = 41
$ws.Range("G9").Value = 2.5
$ws.Range("H9").Value = 2.88
$ws.Range("I9").Value = 3.2
$ws.Range("J9").Value = 1.11
$ws.Range("K9").Value = 6.5
$ws.Range("L9").Value = 1.53
$ws.Range("M9").Value = 2.38
$ws.Range("N9").Value = 2.7
$ws.Range("O9").Value = 1.44
$ws.Range("R9").Value = 2.2
$ws.Range("S9").Value = 1.62
$ws.Range("U9").Value = 10
$ws.Range("W9").Value = 23
$ws.Range("Z9").Value = 6

# Row 11
$ws.Range("AH11").Value = 51
$ws.Range("G11").Value = 1.83
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 4.2
$ws.Range("J11").Value = 1.08
$ws.Range("K11").Value = 8

# Row 12
$ws.Range("N12").Value = 1.82
$ws.Range("O12").Value = 1.92

# Row 13
$ws.Range("J13").Value = 1.07
$ws.Range("K13").Value = 9

# Row 14
$ws.Range("N14").Value = 2
$ws.Range("O14").Value = 1.77

# Row 15
$ws.Range("AA15").Value = 6.5
$ws.Range("AB15").Value = 13
$ws.Range("AC15").Value = 41
$ws.Range("AD15").Value = 201
$ws.Range("AE15").Value = 8
$ws.Range("AF15").Value = 12
$ws.Range("AI15").Value = 19
$ws.Range("G15").Value = 2.9
$ws.Range("H15").Value = 3.3
$ws.Range("J15").Value = 1.05
$ws.Range("K15").Value = 11
$ws.Range("L15").Value = 1.29
$ws.Range("M15").Value = 3.5
$ws.Range("N15").Value = 1.98
$ws.Range("O15").Value = 1.88
$ws.Range("P15").Value = 1.4
$ws.Range("Q15").Value = 2.75
$ws.Range("R15").Value = 1.73
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 9.5
$ws.Range("W15").Value = 29
$ws.Range("X15").Value = 23
$ws.Range("Z15").Value = 10

# Row 16
$ws.Range("AE16").Value = 6.9
$ws.Range("AF16").Value = 5.2
$ws.Range("AG16").Value = 9
$ws.Range("AH16").Value = 5.4
$ws.Range("G16").Value = 15.5
$ws.Range("H16").Value = 6.1
$ws.Range("N16").Value = 1.44
$ws.Range("O16").Value = 2.6
$ws.Range("R16").Value = 2.33
$ws.Range("S16").Value = 1.54
$ws.Range("T16").Value = 35
$ws.Range("U16").Value = 120
$ws.Range("V16").Value = 45
$ws.Range("W16").Value = 600
$ws.Range("X16").Value = 200

# Row 17
$ws.Range("AA17").Value = 5.3
$ws.Range("AB17").Value = 12.5
$ws.Range("AC17").Value = 55
$ws.Range("AD17").Value = 400
$ws.Range("AE17").Value = 7.2
$ws.Range("AG17").Value = 9.25
$ws.Range("AH17").Value = 29
$ws.Range("AI17").Value = 22
$ws.Range("AJ17").Value = 29
$ws.Range("G17").Value = 2.2
$ws.Range("H17").Value = 3.1
$ws.Range("I17").Value = 3
$ws.Range("L17").Value = 1.35
$ws.Range("M17").Value = 2.95
$ws.Range("N17").Value = 2.05
$ws.Range("O17").Value = 1.62
$ws.Range("P17").Value = 1.4
$ws.Range("Q17").Value = 2.42
$ws.Range("R17").Value = 1.86
$ws.Range("S17").Value = 1.85
$ws.Range("T17").Value = 5.8
$ws.Range("U17").Value = 8.5
$ws.Range("V17").Value = 7.7
$ws.Range("X17").Value = 15.5
$ws.Range("Y17").Value = 26
$ws.Range("Z17").Value = 8

# Row 18
$ws.Range("AA18").Value = 6
$ws.Range("AF18").Value = 19
$ws.Range("G18").Value = 1.95
$ws.Range("H18").Value = 3.25
$ws.Range("I18").Value = 4.1
$ws.Range("R18").Value = 1.91
$ws.Range("S18").Value = 1.91
$ws.Range("U18").Value = 9
$ws.Range("V18").Value = 9
$ws.Range("W18").Value = 17

# Row 19
$ws.Range("AG19").Value = 19
$ws.Range("I19").Value = 7
$ws.Range("N19").Value = 1.73
$ws.Range("O19").Value = 2.08
$ws.Range("R19").Value = 1.95
$ws.Range("S19").Value = 1.8
$ws.Range("U19").Value = 7

# Row 20
$ws.Range("AD20").Value = 151
$ws.Range("AE20").Value = 9.5
$ws.Range("AI20").Value = 19
$ws.Range("AJ20").Value = 26
$ws.Range("H20").Value = 3.5
$ws.Range("N20").Value = 1.85
$ws.Range("O20").Value = 1.95
$ws.Range("P20").Value = 1.36
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = 1.67
$ws.Range("S20").Value = 2.1
$ws.Range("Y20").Value = 26
$ws.Range("Z20").Value = 12

# Row 21
$ws.Range("AD21").Value = 1250
$ws.Range("AE21").Value = 6.5
$ws.Range("AF21").Value = 10
$ws.Range("G21").Value = 3.3
$ws.Range("I21").Value = 2.3
$ws.Range("W21").Value = 34
$ws.Range("X21").Value = 29
